$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap rows 47 and 48 (Almeria <-> Lugo), including their "Casos activos" values.
# Before: A47=Almeria, C47=72 ; A48=Lugo, C48=5
# After:  A47=Lugo,    C47=5  ; A48=Almeria, C48=72
$ws.Range("A47").Value = "Lugo"
$ws.Range("C47").Value = 5
$ws.Range("A48").Value = "Almeria"
$ws.Range("C48").Value = 72

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 02:46"
